# [EI-979] Rename "Then_Goto" / "Else_Goto" headers to "Then_Question" / "Else_Question"
# in the Survey.xlsx data dictionary (Sheet1, header row 1, columns I and J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value = "Then_Question"
$ws.Range("J1").Value = "Else_Question"

$ws.Range("L10").Select()
